# devedores.xlsx — apply diff: fix invalid/erroneous entry formats, add new debtor rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: correct an invalid entry (20.0 -> 0.0) and append a missing pair ---
$ws.Range("A2").Value = "'0.0"
$ws.Range("O2").Value = "23/12/2020"
$ws.Range("P2").Value = "'-20"

# --- Row 4: new debtor "Kiko" ---
$ws.Range("A4").Value = "'250"
$ws.Range("B4").Value = "Kiko"
$ws.Range("C4").Value = "23/12/2020"
$ws.Range("D4").Value = "'250"

# --- Row 5: new debtor "Ana" (with a real numeric value in D5) ---
$ws.Range("A5").Value = "'220.0"
$ws.Range("B5").Value = "Ana"
$ws.Range("C5").Value = "23/12/2020"
$ws.Range("D5").Value = 200
$ws.Range("E5").Value = "23/12/2020"
$ws.Range("F5").Value = "'70.0"
$ws.Range("G5").Value = "23/12/2020"
$ws.Range("H5").Value = "'70.0"
$ws.Range("I5").Value = "23/12/2020"
$ws.Range("J5").Value = "'70.0"
$ws.Range("K5").Value = "23/12/2020"
$ws.Range("L5").Value = "'70.0"
$ws.Range("M5").Value = "23/12/2020"
$ws.Range("N5").Value = "'70.0"
$ws.Range("O5").Value = "23/12/2020"
$ws.Range("P5").Value = "'70.0"
$ws.Range("Q5").Value = "23/12/2020"
$ws.Range("R5").Value = "'70.0"
$ws.Range("S5").Value = "23/12/2020"
$ws.Range("T5").Value = "'-50.0"
